$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update instance size for RHEL 7 row (row 3) from t2.nano to t2.micro
$ws.Range("C3").Value = "t2.micro"

# Update the active selection to C4, as recorded in the sheet view
$ws.Range("C4").Select()
